$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 20 (shifts MWG1/WAT1/WAT2 etc. down by 3),
# to make room for new scenes ROT7, ROT8, ROT9.
$ws.Rows("20:22").Insert()

# Row 20: ROT7
$ws.Cells.Item(20, 1).Value = "ROT7"
$ws.Cells.Item(20, 2).Value = "location:current:swamp"
$ws.Cells.Item(20, 3).Value = "You’re wading through murky water in {location:current:namewiththe}. Suddenly, the water around you begins to move. You’re surrounded by crocodiles, each one as long as your outstretched arms. What do you do?"
$ws.Cells.Item(20, 4).Value = "Fight them off"
$ws.Cells.Item(20, 5).Value = "Run for it"
$ws.Cells.Item(20, 6).Value = "One lunges for you, mouth open. You knock it on the snout with a nearby branch and poke it in the eyes. You scramble madly to shore, waving your branch at any of them that come close to you. One of them bites you in the leg, but a forceful knock between the eyes loosens their grip.
Your leg is bleeding as you stumble out on to rocky ground. The crocodiles don’t seem to want to follow you there.
After putting some distance between yourself and the crocodiles, you wrap up your leg and continue on your way to {location:goal:namewiththe}."
$ws.Cells.Item(20, 7).Value = "You start moving for shore immediately, quickly but without trying to attract attention. One of them starts accelerating toward you as you approach land, its eyes floating above the water as it swims. You scramble madly out of the water and onto rocky ground just before it gets to you. It doesn’t seem to want to follow you there.
After putting some distance between yourself and the crocodiles, you wrap up your leg and continue on your way to {location:goal:namewiththe}."

# Row 21: ROT8
$ws.Cells.Item(21, 1).Value = "ROT8"
$ws.Cells.Item(21, 2).Value = "location:current:plains"
$ws.Cells.Item(21, 3).Value = "In the distance, you see a copse of trees, and you head for them, eager to find some shade from the hot sun. You find shade and a small freshwater pond, an ideal resting place.
Suddenly, an angry badger emerges from a nearby den. It bares it teeth and growls as it approaches, defending its territory against this new intruder."
$ws.Cells.Item(21, 4).Value = "Leave the badger’s territory"
$ws.Cells.Item(21, 5).Value = "Try to pacify it with food"
$ws.Cells.Item(21, 6).Value = "You back away slowly, returning to the hot sun in {location:current:namewiththe}. Glancing over your shoulder, you can see the badger watching you intently from afar as it disappears below the horizon."
$ws.Cells.Item(21, 7).Value = "Taking your food out of your pack, you toss some on the ground near the badger. Initially, it thinks you’re attacking it, and becomes more aggressive. You back off. Eventually, it realizes you offered it food. You toss it some more. Slowly, it becomes less hostile, and lets you share its copse, but much of your food is gone.
The shade is a welcome respite from the hot sun. You stay the rest of the day under the watchful eye of the badger. You even manage to catch a few fish in the pond and collect some firewood. Later that day, with night falling, you set out again for {location:goal:namewiththe}, the badger’s gaze disappearing eventually below the horizon."

# Row 22: ROT9
$ws.Cells.Item(22, 1).Value = "ROT9"
$ws.Cells.Item(22, 2).Value = "location:current:mountain"
$ws.Cells.Item(22, 3).Value = "This high up in {location;current:namewiththe}, the air is thin and cold. Your arms and legs begin to feel weak, and breathing is harder."
$ws.Cells.Item(22, 4).Value = "Take a break to adjust to the altitude"
$ws.Cells.Item(22, 5).Value = "Power through to the other side of the mountain"
$ws.Cells.Item(22, 6).Value = "You sit on a nearby stone for a few hours, resting up. After a while, you feel fresh enough to start again. Your legs still feel weaker than usual, but with your newfound strength, you can push on easily toward {location:goal:namewiththe}."
$ws.Cells.Item(22, 7).Value = "Your legs are weak, but through sheer force of will you make them move, over stones and boulders, through trees and brush. You feel like you’ll drop dead from the exertion, but you press on.
Finally, you see the ground begin to slope down under you. You’ve reached the other side! Your legs regain strength as you hike down {location:current:namewiththe}, toward {location:goal:namewiththe}."

# Row heights to match target layout
$ws.Rows(20).RowHeight = 113.4
$ws.Rows(21).RowHeight = 124.6
$ws.Rows(22).RowHeight = 79.85

# Restore the view selection to match the saved state in the target file
$ws.Range("G22").Select()

